# modified text parsing code
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the erroneous duplicated "ERROR" values in column C for rows 8 and 9
$ws.Range("C8").Value = ""
$ws.Range("C9").Value = ""

# Fix parsed name: remove stray middle initial "S" from "ISAAC S SANCHEZ"
$ws.Range("C12").Value = "ISAAC SANCHEZ"
$ws.Range("C22").Value = "ISAAC SANCHEZ"
